# This workbook holds a weekly price report for "Durazno" (peach) at the
# "Macroferia Regional de Talca" market. Two new report rows were added
# (for "Florida King" Especial / Primera, dated 2021-11-24), pushing all
# the subsequent rows (previously 145-192) down by two positions
# (becoming 147-194).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before the current row 145; this shifts the old
# rows 145..192 down to 147..194, carrying all of their existing values,
# formatting and styles along with them.
$ws.Range("145:146").Insert()

# Row 145: Durazno, Florida King, Especial
$ws.Cells.Item(145, 4).Value  = 44524
$ws.Cells.Item(145, 5).Value  = 7
$ws.Cells.Item(145, 6).Value  = "Fruta"
$ws.Cells.Item(145, 7).Value  = 100103
$ws.Cells.Item(145, 8).Value  = "Frutos de hueso (carozo)"
$ws.Cells.Item(145, 9).Value  = 100103004
$ws.Cells.Item(145, 10).Value = "Durazno"
$ws.Cells.Item(145, 11).Value = "Florida King"
$ws.Cells.Item(145, 12).Value = "Especial"
$ws.Cells.Item(145, 1).Value  = 5
$ws.Cells.Item(145, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(145, 3).Value  = "Maule"
$ws.Cells.Item(145, 13).Value = 100
$ws.Cells.Item(145, 14).Value = 20000
$ws.Cells.Item(145, 15).Value = 20000
$ws.Cells.Item(145, 16).Value = 20000
$ws.Cells.Item(145, 17).Value = "$/bandeja 15 kilos empedrada"
$ws.Cells.Item(145, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(145, 19).Value = 1333
$ws.Cells.Item(145, 20).Value = 15

# Row 146: Durazno, Florida King, Primera
$ws.Cells.Item(146, 1).Value  = 5
$ws.Cells.Item(146, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(146, 3).Value  = "Maule"
$ws.Cells.Item(146, 4).Value  = 44524
$ws.Cells.Item(146, 5).Value  = 7
$ws.Cells.Item(146, 6).Value  = "Fruta"
$ws.Cells.Item(146, 7).Value  = 100103
$ws.Cells.Item(146, 8).Value  = "Frutos de hueso (carozo)"
$ws.Cells.Item(146, 9).Value  = 100103004
$ws.Cells.Item(146, 10).Value = "Durazno"
$ws.Cells.Item(146, 11).Value = "Florida King"
$ws.Cells.Item(146, 12).Value = "Primera"
$ws.Cells.Item(146, 13).Value = 60
$ws.Cells.Item(146, 14).Value = 18000
$ws.Cells.Item(146, 15).Value = 18000
$ws.Cells.Item(146, 16).Value = 18000
$ws.Cells.Item(146, 17).Value = "$/bandeja 15 kilos empedrada"
$ws.Cells.Item(146, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(146, 19).Value = 1200
$ws.Cells.Item(146, 20).Value = 15

# Column D (Fecha) is a date column; make sure both new rows keep the same
# date-time number format used by the rest of the column.
$ws.Range("D145:D146").NumberFormat = $ws.Range("D147").NumberFormat
